$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Developers Close On `$125 Million Loan For 70–28 Grand Central Parkway in Forest Hills, Queens"
$ws.Range("B2").Value = "https://newyorkyimby.com/2026/01/developers-close-on-125-million-loan-for-70-28-grand-central-parkway-in-forest-hills-queens.html"
$ws.Range("C2").Value = "Naftali Credit Partners and Madison Realty Capital have closed on a `$125 million loan to develop 70–28 Grand Central Parkway, a 13-story residential building in <a href=`"https://newyorkyimby.com/neighborhoods/forest-hills`">Forest Hills</a>, Queens. Designed by SLCE Architects, the 252,399-square-foot structure will yield 241 condominium units in studio to three-bedroom layouts. The project will also include a 100-vehicle parking garage and a collection of amenities. The 32,000-square-foot development site is located to the west of Grand Central Parkway between Jewel Avenue and 72nd Road."
$ws.Range("D2").Value = "2026-01-10T13:00:24+00:00"
$ws.Range("E2").Value = "Sat, 10 Jan 2026 13:00:24 +0000"
$ws.Range("F2").Value = "YIMBY"
$ws.Range("G2").Value = "YIMBY - Forest Hills"
$ws.Range("H2").Value = ""
